$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (ECs -> Pomc -> Mc4r -> ECs) ---
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pomc"
$ws.Range("C2").Value = "Mc4r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.703265666666667
$ws.Range("H2").Value = 5.109797
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.006923666666666667
$ws.Range("N2").Value = 0.020771
$ws.Range("O2").Value = 0.01563438526027703
$ws.Range("P2").Value = 0.01563438526027703
$ws.Range("Q2").Value = 0.01179284372077778
$ws.Range("R2").Value = 0.106135593487
$ws.Range("S2").Value = 0.01563438526027703
$ws.Range("T2").Value = 0.01563438526027703

# --- Update row 3 (ECs -> Pomc -> Mc4r -> MuSCs) ---
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pomc"
$ws.Range("C3").Value = "Mc4r"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.703265666666667
$ws.Range("H3").Value = 5.109797
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.435925
$ws.Range("N3").Value = 1.307775
$ws.Range("O3").Value = 0.9843656147397229
$ws.Range("P3").Value = 0.9843656147397229
$ws.Range("Q3").Value = 0.7424960857416666
$ws.Range("R3").Value = 6.682464771675
$ws.Range("S3").Value = 0.9843656147397229
$ws.Range("T3").Value = 0.9843656147397229

# --- Remove rows 4 and 5 (MuSCs as sending cluster) entirely ---
$ws.Range("A4:A5").EntireRow.Delete()
